$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    1.036527202731309,
    1.038842205395358,
    1.326325506664252,
    1.779657095697206,
    1.809475044738209,
    1.663556186815498,
    1.702194396450396,
    1.947660774265794,
    1.983895126755568,
    1.850913905971616,
    1.590063071725434,
    1.336797387636375,
    1.194900433558507,
    1.059998552217474,
    0.9700298836639007,
    0.9068655966232724,
    0.8513422548435227,
    0.8252901890807021,
    0.7929476054354004,
    0.7749617135886168,
    0.7596660551172044,
    0.7341157469545613,
    0.7181107202104252,
    0.7103191777090804,
    0.7038173918389828,
    0.6989365637979807,
    0.6932520569524969,
    0.6910809928569983,
    0.6875199864733751,
    0.6861444759311689,
    0.6842692588668601,
    0.68330026243507,
    0.6824158713357489,
    0.6819083092044356,
    0.6816123456854286,
    0.6815023151766866,
    0.6811745062053597,
    0.6809586623367709,
    0.6808763781609285,
    0.6808466284810009,
    0.6807751233375516,
    0.680758020328615,
    0.6807259325172021,
    0.6807052985319758,
    0.6806873419142903,
    0.6806737071097959,
    0.6806676814022236,
    0.6806608703856351,
    0.6806589824864423,
    0.6806575008705177
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
